$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 146
$ws.Range("A1:B25").Font.Bold = $false
$ws.Range("B12").Select()
